$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 20

# Row 3
$ws.Range("D3").Value = 44320
$ws.Range("M3").Value = 30

# Row 4
$ws.Range("D4").Value = 44322
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("S4").Value = 1714

# Row 5
$ws.Range("D5").Value = 44322
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("S5").Value = 1143

# Row 6
$ws.Range("D6").Value = 44302
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 2143

# Row 7
$ws.Range("D7").Value = 44302
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 1714

# Row 8
$ws.Range("D8").Value = 44300
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 2143

# Row 9
$ws.Range("D9").Value = 44300

# Row 12
$ws.Range("D12").Value = 44299
$ws.Range("M12").Value = 80
$ws.Range("R12").Value = "Provincia de Santiago"

# Row 13
$ws.Range("D13").Value = 44299
$ws.Range("M13").Value = 75
$ws.Range("R13").Value = "Provincia de Santiago"

# Row 14
$ws.Range("D14").Value = 44301
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 14000
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 2000

# Row 15
$ws.Range("D15").Value = 44301
$ws.Range("M15").Value = 80
$ws.Range("R15").Value = "Región Metropolitana"
